$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet had an AutoFilter applied (colId 2 and colId 6) which hid most
# data rows. Clear the filter criteria so every row is shown again and the
# previously-hidden rows lose their "hidden" flag, while the AutoFilter
# range itself (A1:L87) stays in place.
$ws.ShowAllData()

# Rows 90:149 held a leftover dynamic-array COUNTIF formula (B90) plus a
# long tail of stray zero cells below the real data (which ends at row 87).
# Remove that whole block so the sheet's used range shrinks back down.
$ws.Rows("90:149").Delete()

# Move the selection/cursor to C5, matching the new view state.
$ws.Range("C5").Select()
